# Upload new version with timestamp
# Inserts a new inventory-item row ("مناديل جيب مبلله" / wet wipes tissues) as
# the 4th data row (row 10), pushing the existing "total" row and the footer
# row down by one. Updates the total and refreshes the generated-on timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Push everything from row 10 down by one row, making room for the new
#    data row. (This automatically shifts the old total row 10->11 and the
#    footer row 11->12, merged ranges included.)
$ws.Rows.Item(10).Insert()

# 2) Seed the new row 10 with the same look & feel (styles/borders/fill) as
#    the data row directly above it (row 9), then overwrite the values.
$ws.Range("A9:Q9").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A10").Value2 = 4
$ws.Range("C10").Value2 = "مناديل جيب مبلله "
$ws.Range("H10").Value2 = "8:0"

# Column L holds its numbers as text (matches existing rows 7-9) -- force
# text storage, write the value, then restore the numeric-look style from
# the row above so the cell format matches its neighbours.
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value2 = "0"
$ws.Range("L9").Copy()
$ws.Range("L10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("N10").Value2 = "6.00"

# Column P also stores its number as text; same trick as column L.
$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value2 = "6.0000"
$ws.Range("P9").Copy()
$ws.Range("P10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column Q reuses the same "1:0" transactions-count text as row 9.
$ws.Range("Q10").NumberFormat = "@"
$ws.Range("Q10").Value2 = "1:0"
$ws.Range("Q9").Copy()
$ws.Range("Q10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Re-merge the cell groups for the new data row (same pattern used by
#    rows 7-9).
$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

# 4) Fix up row heights: the new data row uses the height the old total row
#    used to have, and the (now shifted) total row uses the height the data
#    row above the new row used to have. The footer row keeps its height.
$ws.Rows.Item(10).RowHeight = 24.75
$ws.Rows.Item(11).RowHeight = 25.5
$ws.Rows.Item(12).RowHeight = 16.5

# 5) Update the (now shifted) total in P11 to include the new item.
$ws.Range("P11").Value2 = 187.13999999999999

# 6) Refresh the "generated on" timestamp in the footer row (now row 12,
#    since the footer shifted down when the new data row was inserted).
$ws.Range("A12").Value2 = "Sunday, 17 August, 2025 9:44 AM"
